$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 232:233, shifting existing rows 232-338 down to 234-340
$ws.Range("A232:A233").EntireRow.Insert()

# Fill in the new row 232 (Ajo / Chino / Extra, fecha 2021-09-13)
$ws.Range("A232").Value = 6
$ws.Range("B232").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C232").Value = "Metropolitana"
$ws.Range("D232").Value = 44452
$ws.Range("D232").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E232").Value = 13
$ws.Range("F232").Value = 100112003
$ws.Range("G232").Value = "Ajo"
$ws.Range("H232").Value = "Chino"
$ws.Range("I232").Value = "Extra"
$ws.Range("J232").Value = 1100
$ws.Range("K232").Value = 16000
$ws.Range("L232").Value = 17000
$ws.Range("M232").Value = 16545
$ws.Range("N232").Value = "$/malla 10 kilos"
$ws.Range("O232").Value = "China"
$ws.Range("P232").Value = 1654
$ws.Range("Q232").Value = 10
$ws.Range("R232").Value = "Hortaliza"

# Fill in the new row 233 (Ajo / Chino / Primera, fecha 2021-09-13)
$ws.Range("A233").Value = 6
$ws.Range("B233").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C233").Value = "Metropolitana"
$ws.Range("D233").Value = 44452
$ws.Range("D233").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E233").Value = 13
$ws.Range("F233").Value = 100112003
$ws.Range("G233").Value = "Ajo"
$ws.Range("H233").Value = "Chino"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 2800
$ws.Range("K233").Value = 14000
$ws.Range("L233").Value = 14500
$ws.Range("M233").Value = 14268
$ws.Range("N233").Value = "$/caja 10 kilos"
$ws.Range("O233").Value = "China"
$ws.Range("P233").Value = 1427
$ws.Range("Q233").Value = 10
$ws.Range("R233").Value = "Hortaliza"
